$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.56000000000056
$ws.Range("G2").Value = 0.0007314638132400741
$ws.Range("H2").Value = 0.002059400479961557
$ws.Range("K2").Value = 5.665006736508652
$ws.Range("L2").Value = "[2.232916265965928, 9.097097207051377]"
$ws.Range("M2").Value = 0.00133932563570327
$ws.Range("N2").Value = 0.00133932563570327
$ws.Range("O2").Value = -1.698158191129618
$ws.Range("P2").Value = "[-2.465474114528927, -0.9308422677303092]"
$q2 = 0.00002085172319210571
$ws.Range("Q2").Value = $q2
$ws.Range("R2").Value = $q2
$ws.Range("S2").Value = 11.32773320727654
$ws.Range("T2").Value = "[9.301928832088334, 13.353537582464748]"
$ws.Range("W2").Value = 6.908108108108259
$ws.Range("X2").Value = 3.786666666666749
$ws.Range("Y2").Value = 10.02954954954977

# Row 3 updates
$ws.Range("G3").Value = 0.00003432038754369593
$ws.Range("H3").Value = 0.0004763612275136466
$ws.Range("K3").Value = 6.292756536875356
$ws.Range("L3").Value = "[2.739442004605296, 9.846071069145415]"
$ws.Range("M3").Value = 0.0005782535037424097
$ws.Range("N3").Value = 0.001156507007484819
$ws.Range("O3").Value = 2.547237286694427
$ws.Range("P3").Value = "[1.9811845563178876, 3.113290017070966]"
$ws.Range("Q3").Value = 0.0000000000000002220446049250313
$ws.Range("R3").Value = 0.0000000000000004440892098500626
$ws.Range("S3").Value = 11.61354373175633
$ws.Range("T3").Value = "[9.738057169004623, 13.48903029450803]"
$ws.Range("W3").Value = 13.66972972972982
$ws.Range("X3").Value = 11.59855855855864
$ws.Range("Y3").Value = 15.74090090090101
